# QA overhaul: fix HTML template, diversify scripts, improve content quality
# Applies targeted text updates across the NRWaves, cumcontrol, and boosters sheets.

$wb = $excel.ActiveWorkbook

# ---- NRWaves sheet ----
$ws = $wb.Worksheets.Item("NRWaves")
$ws.Range("B2").Value = "still thinking about our conversation... come back when you can"
$ws.Range("B3").Value = "hey, just checking in on you 😊"
$ws.Range("B4").Value = "I'm starting to think you forgot about me..."
$ws.Range("B5").Value = "literally just took something crazy and you're not even here to see it"
$ws.Range("B6").Value = "hey you 😊"

# ---- cumcontrol sheet ----
$ws = $wb.Worksheets.Item("cumcontrol")

# Row 2: delay2 -> reveal2
$ws.Range("A2").Value = "reveal2"
$ws.Range("B2").Value = "patience babe... I'm not showing you everything at once 😊 the best part is coming"
$ws.Range("C2").Value = "REVEAL variant."

# Row 3: delay1 -> reveal1
$ws.Range("A3").Value = "reveal1"
$ws.Range("B3").Value = "hold on... the next one is worth waiting for, trust me"
$ws.Range("C3").Value = "REVEAL. Send next PPV."

# Row 4: sync2 -> buildup2
$ws.Range("A4").Value = "buildup2"
$ws.Range("B4").Value = "I'm taking my time... good things come to those who wait babe"
$ws.Range("C4").Value = "BUILDUP variant."

# Row 5: sync1 -> buildup1
$ws.Range("A5").Value = "buildup1"
$ws.Range("B5").Value = "you want to see more? then you have to wait for it... 😊"
$ws.Range("C5").Value = "BUILDUP. Final PPV."

# Row 6: edge2 -> tease2
$ws.Range("A6").Value = "tease2"
$ws.Range("B6").Value = "don't rush... I want you to enjoy every single moment of this"
$ws.Range("C6").Value = "TEASE variant."

# Row 7: edge1 -> tease1
$ws.Range("A7").Value = "tease1"
$ws.Range("B7").Value = "not yet babe... I want to build this up more first 😊"
$ws.Range("C7").Value = "TEASE. More PPVs left. She controls the pace of revealing."

# ---- boosters sheet ----
$ws = $wb.Worksheets.Item("boosters")
$ws.Range("B3").Value = "I need more"
$ws.Range("B6").Value = "I'm losing my mind because of you"
$ws.Range("B7").Value = "oh my god"
